# "Atualiza instrução de trabalho"
# Filters the "ITI" backlog table down to the July/2025 "Backlog" (column H)
# records and moves the active selection, mirroring what a user does in the
# Excel UI via Data > Filter > Date Filters > "July" on that column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ITI")

# Column H ("Backlog") holds dates; filter the existing A1:K37 AutoFilter
# range down to the rows whose Backlog date falls in July 2025 (rows 20 and
# 30 in the source data), which hides all the other data rows.
[void]$ws.Range("A1:K37").AutoFilter(8, "Jul-25", 7)

# Move the selection, matching the post-edit cursor position.
[void]$ws.Range("L45").Select()
